$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 4320.3335
$ws.Range("I6").Value = 925
$ws.Range("K6").Value = 2775
$ws.Range("M6").Value = -2663

$ws.Range("H112").Value = 6061639.5
$ws.Range("J112").Value = 6819259.5
$ws.Range("L112").Value = 20457778.5
$ws.Range("N112").Value = -20459994.5

$ws.Range("H129").Value = 1257.9584
$ws.Range("I129").Value = 399
$ws.Range("J129").Value = 1380.6666
$ws.Range("K129").Value = 1197
$ws.Range("L129").Value = 4141.9998
$ws.Range("M129").Value = 3803
$ws.Range("N129").Value = -14141.9998

$ws.Range("H132").Value = 183636.58
$ws.Range("I132").Value = 190650.7
$ws.Range("J132").Value = 34002
$ws.Range("K132").Value = 571952.1000000001
$ws.Range("L132").Value = 102006
$ws.Range("M132").Value = -569422.1000000001
$ws.Range("N132").Value = -107066

$ws.Range("H141").Value = 1978.525
$ws.Range("I141").Value = 1260.3582
$ws.Range("J141").Value = 5679.846
$ws.Range("K141").Value = 3781.0746
$ws.Range("L141").Value = 17039.538
$ws.Range("M141").Value = 1398.9254
$ws.Range("N141").Value = -27399.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 167133.33
$ws.Range("I5").Value = 200440
$ws.Range("K5").Value = 200440
$ws.Range("M5").Value = -200328

$ws.Range("H32").Value = 14966.949
$ws.Range("I32").Value = 2222.3286
$ws.Range("J32").Value = 114091.78
$ws.Range("K32").Value = 2222.3286
$ws.Range("L32").Value = 114091.78
$ws.Range("M32").Value = -1935.3286
$ws.Range("N32").Value = -114665.78

$ws.Range("H37").Value = 3920

$ws.Range("H45").Value = 1441.7142
$ws.Range("I45").Value = 1238.4
$ws.Range("K45").Value = 1238.4
$ws.Range("M45").Value = -861.4000000000001

$ws.Range("H88").Value = 5468.3
$ws.Range("I88").Value = 3157.2
$ws.Range("K88").Value = 3157.2
$ws.Range("M88").Value = -2751.2

$ws.Range("H91").Value = 5468.3
$ws.Range("I91").Value = 3157.2
$ws.Range("K91").Value = 3157.2
$ws.Range("M91").Value = -1753.2

$ws.Range("H97").Value = 8208.923000000001
$ws.Range("I97").Value = 8208.923000000001
$ws.Range("K97").Value = 8208.923000000001
$ws.Range("M97").Value = -7712.923000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 167133.33
$ws.Range("I4").Value = 200440
$ws.Range("K4").Value = 200440
$ws.Range("M4").Value = -200325

$ws.Range("H86").Value = 5989.28
$ws.Range("I86").Value = 2423
$ws.Range("J86").Value = 11338.7
$ws.Range("K86").Value = 2423
$ws.Range("L86").Value = 11338.7
$ws.Range("M86").Value = -1300
$ws.Range("N86").Value = -13584.7

$ws.Range("H89").Value = 5989.28
$ws.Range("I89").Value = 2423
$ws.Range("J89").Value = 11338.7
$ws.Range("K89").Value = 12115
$ws.Range("L89").Value = 56693.5
$ws.Range("M89").Value = -6499
$ws.Range("N89").Value = -67925.5

$ws.Range("H94").Value = 632.8333
$ws.Range("I94").Value = 447.9394
$ws.Range("J94").Value = 2666.6667
$ws.Range("K94").Value = 447.9394
$ws.Range("L94").Value = 2666.6667
$ws.Range("M94").Value = 3.060600000000022
$ws.Range("N94").Value = -3568.6667

$ws.Range("H134").Value = 52636740
$ws.Range("I134").Value = 142862270
$ws.Range("K134").Value = 428586810
$ws.Range("M134").Value = -428584275

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 49.333332
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 48.8
$ws.Range("K7").Value = 50
$ws.Range("L7").Value = 48.8
$ws.Range("M7").Value = 63
$ws.Range("N7").Value = -274.8

$ws.Range("H19").Value = 300
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -130
$ws.Range("N19").ClearContents()

$ws.Range("H22").Value = 525.8570999999999
$ws.Range("I22").Value = 320.14285
$ws.Range("J22").Value = 731.5714
$ws.Range("K22").Value = 320.14285
$ws.Range("L22").Value = 731.5714
$ws.Range("M22").Value = 29.85714999999999
$ws.Range("N22").Value = -1431.5714

$ws.Range("H24").Value = 300
$ws.Range("I24").Value = 300
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 300
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -130
$ws.Range("N24").ClearContents()

$ws.Range("H39").Value = 900
$ws.Range("I39").Value = 900
$ws.Range("K39").Value = 900
$ws.Range("M39").Value = -509

$ws.Range("H49").Value = 900
$ws.Range("I49").Value = 900
$ws.Range("K49").Value = 900
$ws.Range("M49").Value = -718

$ws.Range("H58").Value = 2003.55
$ws.Range("I58").Value = 1080.9231
$ws.Range("J58").Value = 3717
$ws.Range("K58").Value = 1080.9231
$ws.Range("L58").Value = 3717
$ws.Range("M58").Value = -877.9231
$ws.Range("N58").Value = -4123

$ws.Range("H134").Value = 2061.4443
$ws.Range("I134").Value = 1407.5581
$ws.Range("K134").Value = 4222.6743
$ws.Range("M134").Value = -1687.6743

$ws.Range("H136").Value = 2003.55
$ws.Range("I136").Value = 1080.9231
$ws.Range("J136").Value = 3717
$ws.Range("K136").Value = 3242.7693
$ws.Range("L136").Value = 11151
$ws.Range("M136").Value = -692.7692999999999
$ws.Range("N136").Value = -16251

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76923304
$ws.Range("I2").Value = 82
$ws.Range("K2").Value = 492
$ws.Range("M2").Value = -379

$ws.Range("H13").Value = 170.9
$ws.Range("I13").Value = 101.25
$ws.Range("K13").Value = 303.75
$ws.Range("M13").Value = -135.75

$ws.Range("H23").Value = 827.3125
$ws.Range("I23").Value = 5040.5
$ws.Range("J23").Value = 225.42857
$ws.Range("K23").Value = 15121.5
$ws.Range("L23").Value = 676.28571
$ws.Range("M23").Value = -14886.5
$ws.Range("N23").Value = -1146.28571

$ws.Range("H39").Value = 9078.947
$ws.Range("J39").Value = 9078.947
$ws.Range("L39").Value = 27236.841
$ws.Range("N39").Value = -27824.841

$ws.Range("H55").Value = 3708.25
$ws.Range("J55").Value = 4611
$ws.Range("L55").Value = 13833
$ws.Range("N55").Value = -14187

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 64.2
$ws.Range("I2").Value = 20.5
$ws.Range("J2").Value = 93.333336
$ws.Range("K2").Value = 20.5
$ws.Range("L2").Value = 93.333336
$ws.Range("M2").Value = 92.5
$ws.Range("N2").Value = -319.333336

$ws.Range("H46").Value = 22997.5
$ws.Range("I46").Value = 4990
$ws.Range("J46").Value = 29000
$ws.Range("K46").Value = 4990
$ws.Range("L46").Value = 29000
$ws.Range("M46").Value = -4834
$ws.Range("N46").Value = -29312

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()

$ws.Range("H70").Value = 5786.5
$ws.Range("I70").Value = 5776.8696
$ws.Range("K70").Value = 5776.8696
$ws.Range("M70").Value = -5506.8696

$ws.Range("H73").Value = 5786.5
$ws.Range("I73").Value = 5776.8696
$ws.Range("K73").Value = 5776.8696
$ws.Range("M73").Value = -4840.8696

$ws.Range("H132").Value = 2571.0508
$ws.Range("I132").Value = 2345.7778
$ws.Range("J132").Value = 3295.1428
$ws.Range("K132").Value = 7037.3334
$ws.Range("L132").Value = 9885.428400000001
$ws.Range("M132").Value = -4507.3334
$ws.Range("N132").Value = -14945.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3625.6956
$ws.Range("J40").Value = 4216.647
$ws.Range("L40").Value = 4216.647
$ws.Range("N40").Value = -4488.647

$ws.Range("H74").Value = 193464.67
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31996

$ws.Range("H77").Value = 193464.67
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99984

$ws.Range("H136").Value = 5476.6
$ws.Range("I136").Value = 3946.28
$ws.Range("K136").Value = 11838.84
$ws.Range("M136").Value = -9288.84

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10500
$ws.Range("I29").Value = 9750
$ws.Range("J29").Value = 12000
$ws.Range("K29").Value = 9750
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = -9460
$ws.Range("N29").Value = -12580

$ws.Range("H30").Value = 6500
$ws.Range("J30").Value = 6500
$ws.Range("L30").Value = 6500
$ws.Range("N30").Value = -6714

$ws.Range("H55").Value = 1551
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("M55").Value = -723

$ws.Range("H63").Value = 29787.25
$ws.Range("J63").Value = 29787.25
$ws.Range("L63").Value = 29787.25
$ws.Range("N63").Value = -31035.25

$ws.Range("H66").Value = 29787.25
$ws.Range("J66").Value = 29787.25
$ws.Range("L66").Value = 89361.75
$ws.Range("N66").Value = -95601.75

$ws.Range("H123").Value = 25871.41
$ws.Range("J123").Value = 25871.41
$ws.Range("L123").Value = 25871.41
$ws.Range("N123").Value = -35671.41

$ws.Range("H136").Value = 19433.299
$ws.Range("I136").Value = 26327.875
$ws.Range("J136").Value = 3210.7646
$ws.Range("K136").Value = 78983.625
$ws.Range("L136").Value = 9632.293799999999
$ws.Range("M136").Value = -76433.625
$ws.Range("N136").Value = -14732.2938
